$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append Marc's DTR entry as a new row 3: Name, Date, Time In
# (Time Out is left blank, matching the existing John row above it).
$ws.Cells.Item(3, 1).Value = "Marc"

# Format the Date cell as Text first so the literal "2026-02-11" string is
# stored as-is instead of being auto-converted into a date serial number.
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2026-02-11"

$ws.Cells.Item(3, 3).Value = "15:29:19"
